# Weekly fruit/vegetable update: two new daily price records for
# "Ciboulette" (Vega Central Mapocho de Santiago) need to be inserted
# ahead of the existing history, shifting the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at row 229 (pushes old rows 229-248 down to 230-249,
# dimension grows from R248 to R250).
$ws.Rows.Item(229).Insert()
$ws.Rows.Item(229).Insert()

# New row 229: "Primera" quality entry dated 44461
$ws.Range("A229").Value = 9
$ws.Range("B229").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C229").Value = "Metropolitana"
$ws.Range("D229").Value = 44461
$ws.Range("E229").Value = 13
$ws.Range("F229").Value = 100112039
$ws.Range("G229").Value = "Ciboulette"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 160
$ws.Range("K229").Value = 1800
$ws.Range("L229").Value = 2200
$ws.Range("M229").Value = 2000
$ws.Range("N229").Value = "$/docena de atados"
$ws.Range("O229").Value = "Región Metropolitana"
$ws.Range("P229").Value = 667
$ws.Range("Q229").Value = 3
$ws.Range("R229").Value = "Hortaliza"

# New row 230: "Segunda" quality entry, same date 44461
$ws.Range("A230").Value = 9
$ws.Range("B230").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C230").Value = "Metropolitana"
$ws.Range("D230").Value = 44461
$ws.Range("E230").Value = 13
$ws.Range("F230").Value = 100112039
$ws.Range("G230").Value = "Ciboulette"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Segunda"
$ws.Range("J230").Value = 97
$ws.Range("K230").Value = 1300
$ws.Range("L230").Value = 1700
$ws.Range("M230").Value = 1498
$ws.Range("N230").Value = "$/docena de atados"
$ws.Range("O230").Value = "Región Metropolitana"
$ws.Range("P230").Value = 499
$ws.Range("Q230").Value = 3
$ws.Range("R230").Value = "Hortaliza"
